$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$footerFirst   = $sec.Footers.Item(2)   # footer1.xml - Pearson logo, docPr id=3
$footerPrimary = $sec.Footers.Item(1)   # footer2.xml - Pearson logo, docPr id=2
$headerFirst   = $sec.Headers.Item(2)   # header1.xml - BTec logo,   docPr id=1

function Rename-InlinePicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $floatingShape = $inlineShape.ConvertToShape()
    $floatingShape.Name = $newName
    $floatingShape.ConvertToInlineShape() | Out-Null
}

# Pearson Edexcel logo (first-page footer), id="3": image1.png -> image2.png
Rename-InlinePicture $footerFirst.Range "image2.png"

# Pearson Edexcel logo (default footer), id="2": image1.png -> image2.png
Rename-InlinePicture $footerPrimary.Range "image2.png"

# BTec logo (first-page header), id="1": image2.jpg -> image1.jpg
Rename-InlinePicture $headerFirst.Range "image1.jpg"
